# Weekly crime-stat data refresh (new data collected)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (rich-text shared strings) - update the embedded
# "Volume .. Number .." issue number and the reporting week date range.
# ---------------------------------------------------------------------------

# A8: "Volume 32   Number  49" -> "...50"   ("49" begins at char 21, length 2)
$ws.Cells.Item(8, 1).Characters(21, 2).Text = "50"

# C9: "Report Covering the Week  12/1/2025  Through  12/7/2025"
#     -> "...12/8/2025  Through  12/14/2025"
$ws.Cells.Item(9, 3).Characters(27, 9).Text = "12/8/2025"
$ws.Cells.Item(9, 3).Characters(47, 9).Text = "12/14/2025"

# ---------------------------------------------------------------------------
# Data table refresh (rows 15-31). Columns: C/D = week-to-date 2025/2024,
# E = % chg, F/G = 28-day 2025/2024, H = % chg, I/J = YTD 2025/2024,
# K = % chg, L = 28-day % chg, M = YTD % chg, N = 2yr % chg.
# Some cells flip between a number and the "0" / "***.*" text placeholders;
# for those we copy an existing placeholder/number cell (to carry over the
# correct style) before writing the final value.
# ---------------------------------------------------------------------------

function Set-Num($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

# Row 15
Set-Num 15 3 2
Set-Num 15 5 100
Set-Num 15 6 5
Set-Num 15 7 3
Set-Num 15 8 66.666666666666
Set-Num 15 9 41
Set-Num 15 10 41
Set-Num 15 11 0
Set-Num 15 12 32.258064516129
Set-Num 15 13 95.238095238095
Set-Num 15 14 -43.835616438356

# Row 16
Set-Num 16 3 2
Set-Num 16 4 10
Set-Num 16 5 -80
Set-Num 16 6 23
Set-Num 16 7 47
Set-Num 16 8 -51.063829787234
Set-Num 16 9 407
Set-Num 16 10 450
Set-Num 16 11 -9.555555555555
Set-Num 16 12 0.742574257425
Set-Num 16 13 32.573289902280
Set-Num 16 14 -67.282958199356

# Row 17
Set-Num 17 3 6
Set-Num 17 4 19
Set-Num 17 5 -68.421052631578
Set-Num 17 6 52
Set-Num 17 7 57
Set-Num 17 8 -8.771929824561
Set-Num 17 9 752
Set-Num 17 10 717
Set-Num 17 11 4.881450488145
Set-Num 17 12 11.078286558345
Set-Num 17 13 144.155844155844
Set-Num 17 14 -19.31330472103

# Row 18
Set-Num 18 3 3
Set-Num 18 4 3
Set-Num 18 5 0
Set-Num 18 6 11
Set-Num 18 7 18
Set-Num 18 8 -38.888888888888
Set-Num 18 9 208
Set-Num 18 10 279
Set-Num 18 11 -25.448028673835
Set-Num 18 12 5.583756345177
Set-Num 18 13 41.496598639455
Set-Num 18 14 -80.524344569288

# Row 19
Set-Num 19 3 8
Set-Num 19 4 12
Set-Num 19 5 -33.333333333333
Set-Num 19 6 42
Set-Num 19 7 44
Set-Num 19 8 -4.545454545454
Set-Num 19 9 556
Set-Num 19 10 534
Set-Num 19 11 4.119850187265
Set-Num 19 12 11.2
Set-Num 19 13 147.111111111111
Set-Num 19 14 57.954545454545

# Row 20
Set-Num 20 3 5
Set-Num 20 4 3
Set-Num 20 5 66.666666666666
Set-Num 20 6 19
Set-Num 20 7 14
Set-Num 20 8 35.714285714285
Set-Num 20 9 245
Set-Num 20 10 260
Set-Num 20 11 -5.769230769230
Set-Num 20 12 -39.506172839506
Set-Num 20 13 116.814159292035
Set-Num 20 14 -46.389496717724

# Row 21
Set-Num 21 3 26
Set-Num 21 4 48
Set-Num 21 5 -45.833333333333
Set-Num 21 6 152
Set-Num 21 7 183
Set-Num 21 8 -16.939890710382
Set-Num 21 9 2219
Set-Num 21 10 2286
Set-Num 21 11 -2.930883639545
Set-Num 21 12 -0.314465408805
Set-Num 21 13 95.851721094439
Set-Num 21 14 -46.632996632996

# Row 22 - D22 and E22 switch from numbers to the "0"/"***.*" placeholders.
$ws.Cells.Item(14, 4).Copy($ws.Cells.Item(22, 4))   # D22 <- style/value like C14 ("0")
$ws.Cells.Item(14, 5).Copy($ws.Cells.Item(22, 5))   # E22 <- style/value like C14's neighbour ("***.*")
Set-Num 22 6 1
Set-Num 22 8 -66.666666666666

# Row 23
Set-Num 23 3 5
Set-Num 23 4 5
Set-Num 23 5 0
Set-Num 23 7 26
Set-Num 23 8 -46.153846153846
Set-Num 23 9 280
Set-Num 23 10 382
Set-Num 23 11 -26.701570680628
Set-Num 23 12 -30.348258706467
Set-Num 23 13 46.596858638743

# Row 24
Set-Num 24 3 41
Set-Num 24 4 26
Set-Num 24 5 57.692307692307
Set-Num 24 6 114
Set-Num 24 7 104
Set-Num 24 8 9.615384615384
Set-Num 24 9 1125
Set-Num 24 10 1015
Set-Num 24 11 10.837438423645
Set-Num 24 12 -0.968309859154
Set-Num 24 13 50.804289544235

# Row 25
Set-Num 25 4 3
Set-Num 25 5 33.333333333333
Set-Num 25 6 16
Set-Num 25 8 23.076923076923
Set-Num 25 9 204
Set-Num 25 10 165
Set-Num 25 11 23.636363636363
Set-Num 25 12 -20

# Row 26
Set-Num 26 3 18
Set-Num 26 4 16
Set-Num 26 5 12.5
Set-Num 26 6 59
Set-Num 26 7 81
Set-Num 26 8 -27.160493827160
Set-Num 26 9 836
Set-Num 26 10 953
Set-Num 26 11 -12.277019937040
Set-Num 26 12 -15.979899497487
Set-Num 26 13 -1.762632197414

# Row 27
Set-Num 27 3 2
Set-Num 27 5 100
Set-Num 27 6 5
Set-Num 27 7 3
Set-Num 27 8 66.666666666666
Set-Num 27 9 46
Set-Num 27 10 57
Set-Num 27 11 -19.298245614035
Set-Num 27 12 -11.538461538461

# Row 28
Set-Num 28 3 3
Set-Num 28 4 2
Set-Num 28 5 50
Set-Num 28 6 5
Set-Num 28 7 7
Set-Num 28 8 -28.571428571428
Set-Num 28 9 60
Set-Num 28 10 75
Set-Num 28 11 -20
Set-Num 28 12 -34.065934065934

# Row 29 - C29 switches from a number to the "0" placeholder.
$ws.Cells.Item(29, 4).Copy($ws.Cells.Item(29, 3))   # C29 <- style/value like D29 ("0")
Set-Num 29 14 -70.940170940170

# Row 30 - C30 switches from a number to the "0" placeholder.
$ws.Cells.Item(30, 4).Copy($ws.Cells.Item(30, 3))   # C30 <- style/value like D30 ("0")
Set-Num 30 14 -72.321428571428

# Row 31 - F31 switches from the "0" placeholder to a real number.
$ws.Cells.Item(31, 9).Copy($ws.Cells.Item(31, 6))   # F31 <- style like I31 (number format)
Set-Num 31 6 1
Set-Num 31 9 3
Set-Num 31 11 50
Set-Num 31 12 200
